# API: Imports support question_id
#
# 1. "Challenges" sheet: column headers are renamed from human-readable
#    Title Case labels to the snake_case field names used by the API
#    import (id, name, short_description, ... skill_longtermvision), and
#    the skill-weighting columns (K:T) move from large integer "points"
#    values to small fractional weights; "Show Statistics Continuously"
#    becomes a plain numeric 1/0 instead of the string "true".
# 2. "Questions" sheet: headers are likewise renamed to snake_case and a
#    new leading "question_id" column is inserted; the surviving
#    question row's `choices` JSON is cleaned up (no more escaped
#    slashes / explicit `"image":null`); the two extra demo question
#    rows (multi_select / numeric / text) are removed, leaving only the
#    header row and the single_select example row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Challenges"
# ---------------------------------------------------------------------
$challenges = $wb.Worksheets.Item("Challenges")

# Rename headers to the snake_case API field names (column order is
# unchanged; only the label text changes).
$challenges.Range("A1").Value = "id"
$challenges.Range("B1").Value = "name"
$challenges.Range("C1").Value = "short_description"
$challenges.Range("D1").Value = "description"
$challenges.Range("E1").Value = "image"
$challenges.Range("F1").Value = "max_points"
$challenges.Range("G1").Value = "starts_at"
$challenges.Range("H1").Value = "expires_at"
$challenges.Range("I1").Value = "hint_text"
$challenges.Range("J1").Value = "hint_image"
$challenges.Range("K1").Value = "skill_analytical"
$challenges.Range("L1").Value = "skill_strategicplanning"
$challenges.Range("M1").Value = "skill_adaptability"
$challenges.Range("N1").Value = "skill_premierleagueknowledge"
$challenges.Range("O1").Value = "skill_riskmanagement"
$challenges.Range("P1").Value = "skill_decisionmakingunderpressure"
$challenges.Range("Q1").Value = "skill_financialmanagement"
$challenges.Range("R1").Value = "skill_longtermvision"
$challenges.Range("S1").Value = "show_statistics_continuously"
$challenges.Range("T1").Value = "gameweek"

# Rewrite the skill-weight / flag values in row 2.
$challenges.Range("K2").Value = 1       # show_statistics_continuously
$challenges.Range("L2").Value = 1       # gameweek
$challenges.Range("M2").Value = 0.25    # skill_analytical
$challenges.Range("N2").Value = 0.3     # skill_strategicplanning
$challenges.Range("O2").Value = 0.15    # skill_adaptability
$challenges.Range("P2").Value = 0.4     # skill_premierleagueknowledge
$challenges.Range("Q2").Value = 0.2     # skill_riskmanagement
$challenges.Range("R2").Value = 0.25    # skill_decisionmakingunderpressure
$challenges.Range("S2").Value = 0.35    # skill_financialmanagement
$challenges.Range("T2").Value = 0.1     # skill_longtermvision

# ---------------------------------------------------------------------
# Sheet "Questions"
# ---------------------------------------------------------------------
$questions = $wb.Worksheets.Item("Questions")

# Insert a new column at A, pushing challenge_id/text/type/... right.
$questions.Columns("A").Insert()

# Rename headers to the snake_case API field names and add the new
# leading question_id column.
$questions.Range("A1").Value = "question_id"
$questions.Range("B1").Value = "challenge_id"
$questions.Range("C1").Value = "text"
$questions.Range("D1").Value = "type"
$questions.Range("E1").Value = "image"
$questions.Range("F1").Value = "numeric_type_min"
$questions.Range("G1").Value = "numeric_type_max"
$questions.Range("H1").Value = "choices"
$questions.Range("I1").Value = "choices_min_selections"
$questions.Range("J1").Value = "choices_max_selections"
$questions.Range("K1").Value = "correct_text_answer"
$questions.Range("L1").Value = "correct_numeric_answer"
$questions.Range("M1").Value = "correct_selected_choice_text"
$questions.Range("N1").Value = "correct_selected_choice_texts"
$questions.Range("O1").Value = "correct_ordered_choice_texts"

$questions.Range("A2").Value = "01933333-0000-7000-8000-000000000006"

# Clean up the choices JSON for the remaining single_select question
# (drop the escaped slashes and the explicit "image":null entry).
$questions.Range("H2").Value = '[{"text":"Mohamed Salah","description":"Liverpool star with great form","image":"https://example.com/salah.jpg"},{"text":"Kevin De Bruyne","description":"Man City playmaker"},{"text":"Bruno Fernandes","description":"Manchester United captain"}]'

# Remove the other demo question rows (multi_select, numeric, text),
# keeping only the header and the single_select example.
$questions.Rows("3:5").Delete()
